# Update "141" schedule workbook: 31/12 11:40 scrape (LP1912 + LP1912-215 + 6203-6173)
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "LP1912": 17 new rows (760-776), header timestamp/count refresh
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Cells.Item(2, 1).Value2 = "Última actualización: 31/12/2025 08:40:09"
$ws1.Cells.Item(3, 1).Value2 = "Total filas: 775"

$rows1 = @(
    @(760, "08:39:58", "08:43", "10_OLMOS",        4, "LP1912", "31/12/2025"),
    @(761, "08:39:58", "08:51", "16_SANTA ANA",    12, "LP1912", "31/12/2025"),
    @(762, "08:39:58", "09:02", "17X38_ROMERO",    23, "LP1912", "31/12/2025"),
    @(763, "08:39:58", "09:03", "23_HERNANDEZ",    24, "LP1912", "31/12/2025"),
    @(764, "08:39:58", "09:06", "23_HERNANDEZ",    27, "LP1912", "31/12/2025"),
    @(765, "08:39:58", "09:08", "16_SANTA ANA",    29, "LP1912", "31/12/2025"),
    @(766, "08:39:58", "09:14", "11_ETCHEVERRY",   35, "LP1912", "31/12/2025"),
    @(767, "08:39:58", "09:21", "16_SANTA ANA",    42, "LP1912", "31/12/2025"),
    @(768, "08:39:58", "09:26", "215_EL PELIGRO",  47, "LP1912", "31/12/2025"),
    @(769, "08:39:58", "09:35", "23_HERNANDEZ",    56, "LP1912", "31/12/2025"),
    @(770, "08:39:58", "09:44", "14_ABASTO",       65, "LP1912", "31/12/2025"),
    @(771, "08:39:58", "09:51", "15_ABASTO",       72, "LP1912", "31/12/2025"),
    @(772, "08:39:58", "09:53", "10_OLMOS",        74, "LP1912", "31/12/2025"),
    @(773, "08:39:58", "10:02", "215C_EL PATO",    83, "LP1912", "31/12/2025"),
    @(774, "08:39:58", "10:04", "14_ABASTO",       85, "LP1912", "31/12/2025"),
    @(775, "08:39:58", "10:14", "10_OLMOS",        95, "LP1912", "31/12/2025"),
    @(776, "08:39:58", "10:18", "11_ETCHEVERRY",   99, "LP1912", "31/12/2025")
)

foreach ($r in $rows1) {
    $rowNum = $r[0]
    $ws1.Cells.Item($rowNum, 2).Value2 = $r[1]
    $ws1.Cells.Item($rowNum, 3).Value2 = $r[2]
    $ws1.Cells.Item($rowNum, 4).Value2 = $r[3]
    $ws1.Cells.Item($rowNum, 5).Value2 = $r[4]
    $ws1.Cells.Item($rowNum, 6).Value2 = $r[5]
    $ws1.Cells.Item($rowNum, 7).Value2 = $r[6]
}

# ---------------------------------------------------------------------------
# Sheet "LP1912-215": 2 new rows (56-57), header timestamp/count refresh
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Cells.Item(2, 1).Value2 = "Última actualización: 31/12/2025 08:40:09"
$ws2.Cells.Item(3, 1).Value2 = "Total filas: 56"

$rows2 = @(
    @(56, "31/12/2025", "08:39:58", "09:26", "215_EL PELIGRO", 47, "LP1912"),
    @(57, "31/12/2025", "08:39:58", "10:02", "215C_EL PATO",   83, "LP1912")
)

foreach ($r in $rows2) {
    $rowNum = $r[0]
    $ws2.Cells.Item($rowNum, 2).Value2 = $r[1]
    $ws2.Cells.Item($rowNum, 3).Value2 = $r[2]
    $ws2.Cells.Item($rowNum, 4).Value2 = $r[3]
    $ws2.Cells.Item($rowNum, 5).Value2 = $r[4]
    $ws2.Cells.Item($rowNum, 6).Value2 = $r[5]
    $ws2.Cells.Item($rowNum, 7).Value2 = $r[6]
}

# ---------------------------------------------------------------------------
# Sheet "6203-6173": 1 new row (92), header timestamp refresh (count unchanged)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Cells.Item(2, 1).Value2 = "Última actualización: 31/12/2025 08:40:09"
$ws3.Cells.Item(3, 1).Value2 = "Total filas: 91"

$ws3.Cells.Item(92, 2).Value2 = "31/12/2025"
$ws3.Cells.Item(92, 3).Value2 = "08:40:03"
$ws3.Cells.Item(92, 4).Value2 = "10:09"
$ws3.Cells.Item(92, 5).Value2 = "215C_LA PLATA"
$ws3.Cells.Item(92, 6).Value2 = 89
$ws3.Cells.Item(92, 7).Value2 = "L6203"
